$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text string into a cell without Excel's
# autodetection turning a date-like string (e.g. "03-11-2021") into a
# date serial number / applying a date style. We do this by writing a
# text formula that evaluates to the desired string, then converting
# the formula to a static value with a Values-only paste. This keeps
# the cell as a plain shared-string value with the default style.
function Set-TextCell {
    param($addr, [string]$text)
    $ws.Range($addr).Formula = '="' + $text + '"'
    $ws.Range($addr).Copy()
    $ws.Range($addr).PasteSpecial(-4163)  # xlPasteValues
}

# Row 213 (existing row): add the missing MOVE value and update the VIX value
$ws.Range("B213").Value = 71.04000000000001
$ws.Range("C213").Value = 16.03

# Row 214 (new): 03-11-2021
Set-TextCell "A214" "03-11-2021"
$ws.Range("B214").Value = 71.22
$ws.Range("C214").Value = 15.1

# Row 215 (new): 04-11-2021
Set-TextCell "A215" "04-11-2021"
$ws.Range("B215").Value = 64.61
$ws.Range("C215").Value = 15.44

# Row 216 (new): 05-11-2021
Set-TextCell "A216" "05-11-2021"
$ws.Range("B216").Value = 66.90000000000001
$ws.Range("C216").Value = 16.48

# Row 217 (new): 08-11-2021 - only Serie and VIX, no MOVE value
Set-TextCell "A217" "08-11-2021"
$ws.Range("C217").Value = 16.68

$excel.CutCopyMode = $false
